$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.541.98'
$ws.Range('E2').Value = '  +1.88%  '
$ws.Range('D3').Value = '1.664.78'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '237.83'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4811'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2636'
$ws.Range('E8').Value = '  +0.60%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06169'
$ws.Range('E9').Value = '  +2.47%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07114'
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('D11').Value = '1.661.96'
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.81'
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.5895'
$ws.Range('E13').Value = '  -5.44%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.379'
$ws.Range('E14').Value = '  -4.83%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '75.14'
$ws.Range('E15').Value = '  +2.65%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.9999'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.9995'
$ws.Range('D18').Value = '25.533.18'
$ws.Range('E18').Value = '  +1.88%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000006749'
$ws.Range('E19').Value = '  +1.80%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.48'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').Value = '1.872.00'
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.435'
$ws.Range('E22').Value = '  -2.37%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.717'
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.294'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '134.89'
$ws.Range('E25').Value = '  +2.32%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '15.05'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.377'
$ws.Range('E27').Value = '  -1.52%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '105.22'
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.710'
$ws.Range('E29').Value = '  +1.51%  '
$ws.Range('E30').Value = '  +4.46%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.664'
$ws.Range('E31').Value = '  +1.26%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.07729'
$ws.Range('E32').Value = '  -2.32%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.9990'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04249'
$ws.Range('E34').Value = '  -7.90%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.600'
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6131'
$ws.Range('E36').Value = '  +6.10%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9505'
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.600'
$ws.Range('E38').Value = '  -0.21%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.8637'
$ws.Range('E39').Value = '  +4.89%  '
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.857'
$ws.Range('E41').Value = '  +1.25%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.01466'
$ws.Range('E42').Value = '  -6.18%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '97.04'
$ws.Range('E43').Value = '  -1.50%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.3771'
$ws.Range('E44').Value = '  +0.86%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.853'
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.1124'
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.216'
$ws.Range('E47').Value = '  +1.61%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.05260'
$ws.Range('E48').Value = '  +1.45%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '29.77'
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.373'
$ws.Range('E50').Value = '  +1.93%  '
$ws.Range('E51').Value = '  +0.03%  '
